# Fruta / hortaliza, semanal
# Insert a new weekly record at row 122, shifting the existing rows 122-125 down to 123-126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 122 (this shifts rows 122:125 -> 123:126)
$ws.Rows.Item(122).Insert()

# Populate the new row 122 with the new weekly Ají price record
$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value = "Ñuble"
$ws.Cells.Item(122, 4).Value = 44939
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = 100112021
$ws.Cells.Item(122, 7).Value = "Ají"
$ws.Cells.Item(122, 8).Value = "Americana (o)"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 30
$ws.Cells.Item(122, 11).Value = 13000
$ws.Cells.Item(122, 12).Value = 13000
$ws.Cells.Item(122, 13).Value = 13000
$ws.Cells.Item(122, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(122, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(122, 16).Value = 867
$ws.Cells.Item(122, 17).Value = 15
$ws.Cells.Item(122, 18).Value = "Hortaliza"
